$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 476.8
$ws.Range("I58").Value = 476.8
$ws.Range("K58").Value = 1430.4
$ws.Range("M58").Value = -1280.4
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 2000
$ws.Range("K74").Value = 2000
$ws.Range("M74").Value = -1064
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 2000
$ws.Range("K77").Value = 10000
$ws.Range("M77").Value = -5320
$ws.Range("H80").Value = 2802.2727
$ws.Range("I80").Value = 1031.25
$ws.Range("J80").Value = 3814.2856
$ws.Range("K80").Value = 3093.75
$ws.Range("L80").Value = 11442.8568
$ws.Range("M80").Value = -2095.75
$ws.Range("N80").Value = -13438.8568
$ws.Range("H83").Value = 2802.2727
$ws.Range("I83").Value = 1031.25
$ws.Range("J83").Value = 3814.2856
$ws.Range("K83").Value = 9281.25
$ws.Range("L83").Value = 34328.5704
$ws.Range("M83").Value = -4289.25
$ws.Range("N83").Value = -44312.5704
$ws.Range("H137").Value = 1263.6666
$ws.Range("I137").Value = 1195.5
$ws.Range("J137").Value = 1400
$ws.Range("K137").Value = 3586.5
$ws.Range("L137").Value = 4200
$ws.Range("M137").Value = -1036.5
$ws.Range("N137").Value = -9300
$ws.Range("H138").Value = 1964.7693
$ws.Range("I138").Value = 1362.25
$ws.Range("J138").Value = 2232.5557
$ws.Range("K138").Value = 4086.75
$ws.Range("L138").Value = 6697.6671
$ws.Range("M138").Value = 1053.25
$ws.Range("N138").Value = -16977.6671

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2238.375
$ws.Range("I74").Value = 1903.3334
$ws.Range("J74").Value = 2439.4
$ws.Range("K74").Value = 1903.3334
$ws.Range("L74").Value = 2439.4
$ws.Range("M74").Value = -1029.3334
$ws.Range("N74").Value = -4187.4
$ws.Range("H77").Value = 2238.375
$ws.Range("I77").Value = 1903.3334
$ws.Range("J77").Value = 2439.4
$ws.Range("K77").Value = 9516.666999999999
$ws.Range("L77").Value = 12197
$ws.Range("M77").Value = -5148.666999999999
$ws.Range("N77").Value = -20933
$ws.Range("H97").Value = 1121.0476
$ws.Range("I97").Value = 1019.1429
$ws.Range("J97").Value = 1324.8572
$ws.Range("K97").Value = 1019.1429
$ws.Range("L97").Value = 1324.8572
$ws.Range("M97").Value = -523.1429000000001
$ws.Range("N97").Value = -2316.8572

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2123.4666
$ws.Range("J94").Value = 1135.5714
$ws.Range("L94").Value = 1135.5714
$ws.Range("N94").Value = -2037.5714
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1832.3334
$ws.Range("I31").Value = 1498.5
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1498.5
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1203.5
$ws.Range("N31").Value = -3090
$ws.Range("H34").Value = 1832.3334
$ws.Range("I34").Value = 1498.5
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1498.5
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1296.5
$ws.Range("N34").Value = -2904
$ws.Range("H58").Value = 964.4
$ws.Range("I58").Value = 992.7143
$ws.Range("K58").Value = 992.7143
$ws.Range("M58").Value = -789.7143
$ws.Range("H132").Value = 3497.5
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H136").Value = 964.4
$ws.Range("I136").Value = 992.7143
$ws.Range("K136").Value = 2978.1429
$ws.Range("M136").Value = -428.1428999999998

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 617.75
$ws.Range("J5").Value = 477.5
$ws.Range("L5").Value = 1432.5
$ws.Range("N5").Value = -1656.5
$ws.Range("H39").Value = 4000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 12000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -12588
$ws.Range("H55").Value = 957.1429000000001
$ws.Range("I55").Value = 718.1818
$ws.Range("J55").Value = 1833.3334
$ws.Range("K55").Value = 2154.5454
$ws.Range("L55").Value = 5500.0002
$ws.Range("M55").Value = -1977.5454
$ws.Range("N55").Value = -5854.0002
$ws.Range("H108").Value = 347.83334
$ws.Range("I108").Value = 347.83334
$ws.Range("K108").Value = 1043.50002
$ws.Range("M108").Value = 1836.49998
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("K111").Value = 3000
$ws.Range("M111").Value = 67
$ws.Range("H112").Value = 27796.727
$ws.Range("I112").Value = 817.5
$ws.Range("J112").Value = 43213.43
$ws.Range("K112").Value = 2452.5
$ws.Range("L112").Value = 129640.29
$ws.Range("M112").Value = -1344.5
$ws.Range("N112").Value = -131856.29
$ws.Range("H135").Value = 617.75
$ws.Range("J135").Value = 477.5
$ws.Range("L135").Value = 4297.5
$ws.Range("N135").Value = -9367.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 832.05554
$ws.Range("I2").Value = 1311.125
$ws.Range("J2").Value = 448.8
$ws.Range("K2").Value = 1311.125
$ws.Range("L2").Value = 448.8
$ws.Range("M2").Value = -1198.125
$ws.Range("N2").Value = -674.8
$ws.Range("H122").Value = 24926.732
$ws.Range("I122").Value = 30750.334
$ws.Range("K122").Value = 92251.00199999999
$ws.Range("M122").Value = -89801.00199999999
$ws.Range("H126").Value = 11000
$ws.Range("I126").Value = 7000
$ws.Range("K126").Value = 21000
$ws.Range("M126").Value = -18530

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 444.42856
$ws.Range("J55").Value = 513
$ws.Range("L55").Value = 513
$ws.Range("N55").Value = -859
$ws.Range("H100").Value = 4710.5
$ws.Range("I100").Value = 4710.5
$ws.Range("K100").Value = 4710.5
$ws.Range("M100").Value = -4169.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1612.5
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

Write-Output "Applied all changes"